# Apply the K-column (column G) value updates as described in the diff.
# This corresponds to the commit: "regen save_data to use K instead of Strike#,
# regen std/mean, calc and write s_vals" -- for this particular sheet, the
# observable change is a set of updated numeric values in column G (header "K").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "G2" = 1
    "G4" = 3
    "G5" = 0
    "G6" = 0
    "G8" = 1
    "G9" = 0
    "G10" = 2
    "G12" = 0
    "G13" = 2
    "G14" = 0
    "G15" = 1
    "G16" = 2
    "G17" = 0
    "G18" = 1
    "G19" = 0
    "G20" = 1
    "G21" = 1
    "G22" = 3
    "G23" = 1
    "G24" = 1
    "G25" = 2
    "G26" = 1
    "G27" = 2
    "G28" = 1
    "G29" = 1
    "G30" = 1
    "G31" = 1
    "G32" = 2
    "G33" = 2
    "G34" = 0
    "G35" = 2
    "G36" = 0
    "G37" = 1
    "G39" = 0
    "G40" = 0
    "G41" = 0
    "G43" = 1
    "G44" = 2
    "G46" = 3
    "G47" = 0
    "G48" = 1
    "G49" = 2
    "G50" = 2
    "G51" = 1
    "G52" = 4
    "G53" = 2
    "G54" = 2
    "G55" = 0
    "G56" = 1
    "G57" = 1
    "G58" = 2
    "G59" = 1
    "G60" = 2
    "G62" = 2
    "G63" = 1
    "G64" = 0
    "G65" = 1
    "G66" = 1
    "G67" = 1
    "G68" = 1
    "G69" = 1
    "G70" = 2
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
